$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C4").Value = 0.0434083408340834
$ws.Range("E4").Value = -0.09261326132613261
$ws.Range("F4").Value = 0.04191781626421821
$ws.Range("G4").Value = 0.05524152415241524
$ws.Range("H4").Value = -0.0683108310831083
$ws.Range("I4").Value = -0.01120912091209121
$ws.Range("J4").Value = -0.07935598285662938

$ws.Range("C5").Value = -0.1497269726972697
$ws.Range("E5").Value = 0.0375037503750375
$ws.Range("F5").Value = 0.0084217820368032
$ws.Range("G5").Value = 0.03733573357335733
$ws.Range("H5").Value = -0.06744674467446744
$ws.Range("I5").Value = -0.1615121512151215
$ws.Range("J5").Value = 0.009259398181169021

$ws.Range("C6").Value = 0.05123312331233122
$ws.Range("E6").Value = 0.05489348934893489
$ws.Range("F6").Value = 0.1496265392016802
$ws.Range("G6").Value = 0.04318031803180317
$ws.Range("H6").Value = -0.01503750375037504
$ws.Range("I6").Value = -0.05526552655265525
$ws.Range("J6").Value = 0.072851000595847

$ws.Range("C7").Value = -0.005472547254725472
$ws.Range("E7").Value = 0.01549354935493549
$ws.Range("F7").Value = -0.003166699242363602
$ws.Range("G7").Value = -0.003036303630363036
$ws.Range("H7").Value = 0.005880588058805879
$ws.Range("I7").Value = 0.004116411641164116
$ws.Range("J7").Value = 0.005556839089930339

$ws.Range("C8").Value = 0.06094209420942093
$ws.Range("E8").Value = -0.2776117611761176
$ws.Range("F8").Value = -0.049766316541628
$ws.Range("G8").Value = -0.9696969696969697
$ws.Range("H8").Value = 0.9999759975997599
$ws.Range("I8").Value = 0.7969636963696368
$ws.Range("J8").Value = -0.3627127701000643

$ws.Range("C9").Value = 0.9541554155415541
$ws.Range("E9").Value = 0.06687068706870687
$ws.Range("F9").Value = -0.1074357316234651
$ws.Range("G9").Value = 0.01765376537653765
$ws.Range("H9").Value = 0.000216021602160216
$ws.Range("I9").Value = 0.1024542454245424
$ws.Range("J9").Value = 0.05451823232399259

$ws.Range("C10").Value = -0.09206120612061205
$ws.Range("E10").Value = -0.0808160816081608
$ws.Range("F10").Value = 0.02993349757975595
$ws.Range("G10").Value = -0.09492949294929492
$ws.Range("H10").Value = 0.0553135313531353
$ws.Range("I10").Value = -0.05666966696669666
$ws.Range("J10").Value = -0.094544276308696

$ws.Range("C11").Value = -0.05981398139813981
$ws.Range("E11").Value = 0.1422622262226222
$ws.Range("F11").Value = 0.1059206298307826
$ws.Range("G11").Value = 0.2047644764476447
$ws.Range("H11").Value = -0.2059165916591659
$ws.Range("I11").Value = -0.1155475547554755
$ws.Range("J11").Value = 0.1220404282634053

$ws.Range("C12").Value = -0.04018001800180018
$ws.Range("E12").Value = -0.1270927092709271
$ws.Range("F12").Value = -0.1035592549647097
$ws.Range("G12").Value = 0.05556555655565556
$ws.Range("H12").Value = -0.09106510651065104
$ws.Range("I12").Value = -0.1073267326732673
$ws.Range("J12").Value = -0.1254549438597016

$ws.Range("C13").Value = 0.3081068106810681
$ws.Range("E13").Value = 0.1323612361236123
$ws.Range("F13").Value = -0.03936534747834754
$ws.Range("G13").Value = 0.03152715271527153
$ws.Range("H13").Value = -0.02793879387938793
$ws.Range("I13").Value = 0.564128412841284
$ws.Range("J13").Value = 0.1066120985655534

$ws.Range("C14").Value = -0.1578277827782778
$ws.Range("E14").Value = -0.1305250525052505
$ws.Range("F14").Value = 0.0776250801350078
$ws.Range("G14").Value = -0.03956795679567957
$ws.Range("H14").Value = 0.03041104110411041
$ws.Range("I14").Value = -0.05876987698769877
$ws.Range("J14").Value = -0.1158114876852761

